$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Testing"
$ws.Range("C2").Value = "Testing"
$ws.Range("D2").Value = 13.99

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Testing Second Deal"
$ws.Range("C3").Value = "test Second Deal"
$ws.Range("D3").Value = 25.99

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Testing 3"
$ws.Range("C4").Value = "ertRTGERG"
$ws.Range("D4").Value = 12.99
